$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 values (identifier/title/levelOfDescription/extentAndMedium/notes for MCH213)
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("A2").Value = "MCH213"

$ws.Range("C2:E2").Font.Name = "Calibri"
$ws.Range("C2:E2").Font.Size = 10
$ws.Range("C2:E2").Font.ThemeColor = 1
$ws.Range("C2").Value = "DAS KAIROS DOKUMENT, WAS BUNDESDEUTSCHE BANTEN MIT DER APARTHEID"
$ws.Range("E2").Value = "Series"

$ws.Range("F2").Font.Name = "Calibri"
$ws.Range("F2").Font.Size = 10
$ws.Range("F2").Font.ThemeColor = 1
$ws.Range("F2").Value = "1 Box"

$ws.Range("G2:H2").Font.Name = "Calibri"
$ws.Range("G2:H2").Font.Size = 10
$ws.Range("G2:H2").Font.ThemeColor = 1
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"

# Restore the frozen top row + selection of the new data row
$ws.Range("A2:M2").Select()
$excel.ActiveWindow.FreezePanes = $true
